$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: make room for the new row ---------------------------------
# The new "Macbook / 100 sentences" run is inserted as row 2, pushing the
# existing "europarl / 100000" (old row 2) and "1000000" (old row 3) runs
# down by one row. Inserting at row 3 (instead of row 2) lets the new,
# still-empty row 3 inherit the formatting of row 2 above it (which is
# what the original file already used for the data rows), and the
# following rows shift down automatically.
$ws.Rows.Item(3).Insert()

# --- Step 2: move the old "europarl / 100000" row from row 2 to row 3 --
$ws.Range("A3").Value = "europarl"
$ws.Range("B3").Value = 100000
$ws.Range("C3").Value = 0.19
$ws.Range("D3").Value = 0.91
$ws.Range("E3").Value = 0.55000000000000004
$ws.Range("F3").Value = "7h30"
$ws.Range("G3").Formula = "=480/B3"
$ws.Range("H3").Value = "AWS P3"
$ws.Range("I3").Value = 1000000
$ws.Range("J3").Value = 32215
$ws.Range("K3").Value = 12
$ws.Range("L3").Value = 2
$ws.Range("M3").Value = 4
$ws.Range("N3").Value = 256
$ws.Range("O3").Value = 1024
$ws.Range("P3").Formula = "=((B3*S3)/U3)*T3"
$ws.Range("Q3").Value = "tensor2tensor"
$ws.Range("R3").Value = 16000
$ws.Range("S3").Value = 0.8
$ws.Range("T3").Value = 450
$ws.Range("U3").Value = 80
$ws.Range("V3").Value = 160

# --- Step 3: overwrite row 2 with the new "Macbook / 100 sentences" run
$ws.Range("A2").Value = "europarl"
$ws.Range("B2").Value = 100
$ws.Range("C2").Value = 0.06
$ws.Range("D2").Value = 0.99
$ws.Range("E2").Value = 0.19
$ws.Range("F2").Clear()
$ws.Range("G2").Formula = "=10/B2"
$ws.Range("H2").Value = "Macbook"
$ws.Range("I2").Value = 1000000
$ws.Range("J2").Value = 1013
$ws.Range("K2").Value = 0.06
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2
$ws.Range("N2").Value = 32
$ws.Range("O2").Value = 128
$ws.Range("P2").Formula = "=((B2*S2)/U2)*T2"
$ws.Range("Q2").Value = "tensor2tensor"
$ws.Range("R2").Value = 1000
$ws.Range("S2").Value = 0.8
$ws.Range("T2").Value = 1000
$ws.Range("U2").Value = 10
$ws.Range("V2").Value = 10

# --- Step 4: restore the user's selection on the new data cell ---------
$ws.Range("E2").Select()
